# Lab1/Phisical and Logical Models.pptx — "Add labs for Module 8"
#
# The commit removes the (empty) title slide that used to sit at the
# front of the deck, leaving the six content slides (each holding a
# single full-bleed screenshot) as slides 1-6.

$p = $ppt.ActivePresentation

# Slide 1 is the empty "ctrTitle"/"subTitle" title slide — remove it so
# the picture slides shift up to become slides 1-6.
$p.Slides.Item(1).Delete()
